$d = $word.ActiveDocument

# Find the "LOQ4031: Quimica Geral I (Requisito fraco)" paragraph under
# "Requisitos" -- the three paragraphs that used to follow it (a blank
# paragraph, "Ver no Jupiter Salvar em pdf Salvar em docx" and the
# "(c) 2020 ... Creative Commons Attribution" footer line) were removed
# from the generated page, so delete them from the document too.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4031*Requisito fraco*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $first = $target.Next()              # blank paragraph
    $last  = $first.Next().Next()        # the copyright/footer paragraph

    $killRange = $d.Range($first.Range.Start, $last.Range.End)
    $killRange.Delete()
}
